# Adds a new "team" block (rows 8-13) below the existing one (rows 2-7),
# mirroring its shape: columns A-F stay blank (empty text), column G holds
# the player name.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$names = @(
    "Matteo Mazzola",
    "Stefano Pizzini",
    "Marco Gerola",
    "Michele Parisi",
    "Andrea Anzelini",
    "Davide Raffaelli"
)

$row = 8
foreach ($name in $names) {
    foreach ($col in @("A", "B", "C", "D", "E", "F")) {
        $cell = $ws.Range("$col$row")
        # A bare "'" makes the interpreter store an empty *text* cell
        # (like the existing blank cells) instead of clearing it to an
        # empty/numeric cell; resetting the style afterwards drops the
        # incidental quote-prefix formatting that comes with it.
        $cell.Value = "'"
        $cell.Style = "Normal"
    }
    $ws.Range("G$row").Value = $name
    $row++
}
